# "getting spring event sizes" -- rename the results sheet from
# "storm7_results" to "event 1" and move the live selection/scroll
# position over to where the analysis now looks (N12, instead of the old
# F47/topLeftCell A41).
#
# Note: the chart series references (c:f / xVal / yVal), the bookViews
# window geometry (xWindow/yWindow/windowWidth/windowHeight), the
# xr:revisionPtr documentId and the x15ac:absPath author path are all
# artifacts of the authoring machine / Excel's chart-part + co-authoring
# serializers; they aren't reachable through the Excel object model
# surfaced by this host, so this script sticks to the genuine
# object-model edits: the worksheet name and the selection/view.

$wb = $excel.ActiveWorkbook

$ws = $null
foreach ($sheet in $wb.Worksheets) {
    if ($sheet.Name -eq "storm7_results") {
        $ws = $sheet
    }
}
if ($ws -eq $null) {
    $ws = $wb.Worksheets.Item(1)
}

# Rename the worksheet to match the new "event 1" naming convention.
$ws.Name = "event 1"

# Activate the sheet and move the selection to N12; this also drops the
# stale topLeftCell="A41" scroll anchor that pointed at the old selection.
$ws.Activate() | Out-Null
$ws.Range("N12").Select() | Out-Null
